# Update workbook for data through 2021-11-16 (was through 2021-11-15)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2021-11-16"

# Update the label for the November row
$ws.Range("A12").Value = "November (through 11-16)"

# Update November row (row 12) values
$ws.Range("C12").Value = 39
$ws.Range("D12").Value = 68
$ws.Range("E12").Value = 34
$ws.Range("G12").Value = 93
$ws.Range("H12").Value = 114

# Update Total row (row 13) values
$ws.Range("C13").Value = 525
$ws.Range("D13").Value = 778
$ws.Range("E13").Value = 649
$ws.Range("G13").Value = 1150
$ws.Range("H13").Value = 1556
